$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry arrived for row 80 (Ají, Americana, Limache).
# The previous row-80 entry (Provincia de Limarí) is preserved by
# shifting it down to the newly appended row 81.
$src = $ws.Range("A80:R80")
$dst = $ws.Range("A81:R81")
$src.Copy($dst)

# Now overwrite row 80 with the new entry's values.
$ws.Range("D80").Value = 44628
$ws.Range("J80").Value = 30
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 21000
$ws.Range("M80").Value = 20333
$ws.Range("N80").Value = "$/saco 25 kilos"
$ws.Range("O80").Value = "Limache"
$ws.Range("P80").Value = 813
